$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2095.818
$ws.Cells.Item(19, 9).Value = 580
$ws.Cells.Item(19, 10).Value = 3359
$ws.Cells.Item(19, 11).Value = 580
$ws.Cells.Item(19, 12).Value = 3359
$ws.Cells.Item(19, 13).Value = -405
$ws.Cells.Item(19, 14).Value = -3709

$ws.Cells.Item(80, 8).Value = 1718.3125
$ws.Cells.Item(80, 10).Value = 1367.1666
$ws.Cells.Item(80, 12).Value = 4101.4998
$ws.Cells.Item(80, 14).Value = -6097.4998

$ws.Cells.Item(83, 8).Value = 1718.3125
$ws.Cells.Item(83, 10).Value = 1367.1666
$ws.Cells.Item(83, 12).Value = 12304.4994
$ws.Cells.Item(83, 14).Value = -22288.4994

$ws.Cells.Item(98, 8).Value = 1779.1515
$ws.Cells.Item(98, 9).Value = 1731.5518
$ws.Cells.Item(98, 10).Value = 2124.25
$ws.Cells.Item(98, 11).Value = 1731.5518
$ws.Cells.Item(98, 12).Value = 2124.25
$ws.Cells.Item(98, 13).Value = -233.5518
$ws.Cells.Item(98, 14).Value = -5120.25

$ws.Cells.Item(100, 8).Value = 1246.3
$ws.Cells.Item(100, 9).Value = 751.1429000000001
$ws.Cells.Item(100, 10).Value = 2401.6667
$ws.Cells.Item(100, 11).Value = 751.1429000000001
$ws.Cells.Item(100, 12).Value = 2401.6667
$ws.Cells.Item(100, 13).Value = -210.1429000000001
$ws.Cells.Item(100, 14).Value = -3483.6667

$ws.Cells.Item(111, 8).Value = 7143633.5
$ws.Cells.Item(111, 9).Value = 11111612
$ws.Cells.Item(111, 10).Value = 1272.6
$ws.Cells.Item(111, 11).Value = 33334836
$ws.Cells.Item(111, 12).Value = 3817.8
$ws.Cells.Item(111, 13).Value = -33331769
$ws.Cells.Item(111, 14).Value = -9951.799999999999

$ws.Cells.Item(122, 8).Value = 1779.1515
$ws.Cells.Item(122, 9).Value = 1731.5518
$ws.Cells.Item(122, 10).Value = 2124.25
$ws.Cells.Item(122, 11).Value = 5194.6554
$ws.Cells.Item(122, 12).Value = 6372.75
$ws.Cells.Item(122, 13).Value = -2744.6554
$ws.Cells.Item(122, 14).Value = -11272.75

$ws.Cells.Item(131, 8).Value = 2476.5715
$ws.Cells.Item(131, 9).Value = 777.7273
$ws.Cells.Item(131, 11).Value = 2333.1819
$ws.Cells.Item(131, 13).Value = 2706.8181

$ws.Cells.Item(132, 8).Value = 5377436.5
$ws.Cells.Item(132, 9).Value = 5748137.5
$ws.Cells.Item(132, 10).Value = 2274
$ws.Cells.Item(132, 11).Value = 17244412.5
$ws.Cells.Item(132, 12).Value = 6822
$ws.Cells.Item(132, 13).Value = -17241882.5
$ws.Cells.Item(132, 14).Value = -11882

$ws.Cells.Item(137, 8).Value = 1329.7
$ws.Cells.Item(137, 10).Value = 1350.0454
$ws.Cells.Item(137, 12).Value = 4050.1362
$ws.Cells.Item(137, 14).Value = -9150.136200000001

$ws.Cells.Item(138, 8).Value = 1943.8552
$ws.Cells.Item(138, 9).Value = 1627.1224
$ws.Cells.Item(138, 10).Value = 2518.6667
$ws.Cells.Item(138, 11).Value = 4881.3672
$ws.Cells.Item(138, 12).Value = 7556.000100000001
$ws.Cells.Item(138, 13).Value = 258.6328000000003
$ws.Cells.Item(138, 14).Value = -17836.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2992.55
$ws.Cells.Item(32, 9).Value = 2512.36
$ws.Cells.Item(32, 10).Value = 10195.4
$ws.Cells.Item(32, 11).Value = 2512.36
$ws.Cells.Item(32, 12).Value = 10195.4
$ws.Cells.Item(32, 13).Value = -2225.36
$ws.Cells.Item(32, 14).Value = -10769.4

$ws.Cells.Item(45, 8).Value = 1634.32
$ws.Cells.Item(45, 9).Value = 1496.3125
$ws.Cells.Item(45, 11).Value = 1496.3125
$ws.Cells.Item(45, 13).Value = -1119.3125

$ws.Cells.Item(102, 8).Value = 1183.3334
$ws.Cells.Item(102, 9).Value = 1183.3334
$ws.Cells.Item(102, 11).Value = 1183.3334
$ws.Cells.Item(102, 13).Value = 438.6666

$ws.Cells.Item(110, 8).Value = 3309.0625
$ws.Cells.Item(110, 9).Value = 2351.4285
$ws.Cells.Item(110, 10).Value = 10012.5
$ws.Cells.Item(110, 11).Value = 2351.4285
$ws.Cells.Item(110, 12).Value = 10012.5
$ws.Cells.Item(110, 13).Value = -306.4285
$ws.Cells.Item(110, 14).Value = -14102.5

$ws.Cells.Item(122, 8).Value = 1393.04
$ws.Cells.Item(122, 9).Value = 1092.625
$ws.Cells.Item(122, 10).Value = 2594.7
$ws.Cells.Item(122, 11).Value = 3277.875
$ws.Cells.Item(122, 12).Value = 7784.099999999999
$ws.Cells.Item(122, 13).Value = -827.875
$ws.Cells.Item(122, 14).Value = -12684.1

$ws.Cells.Item(132, 8).Value = 1695.625
$ws.Cells.Item(132, 9).Value = 1313.8334
$ws.Cells.Item(132, 11).Value = 3941.5002
$ws.Cells.Item(132, 13).Value = -1411.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2382.3215
$ws.Cells.Item(105, 9).Value = 2379.4167
$ws.Cells.Item(105, 11).Value = 2379.4167
$ws.Cells.Item(105, 13).Value = -632.4167000000002

$ws.Cells.Item(107, 8).Value = 803.3125
$ws.Cells.Item(107, 9).Value = 741.53845
$ws.Cells.Item(107, 10).Value = 1071
$ws.Cells.Item(107, 11).Value = 741.53845
$ws.Cells.Item(107, 12).Value = 1071
$ws.Cells.Item(107, 13).Value = 1178.46155
$ws.Cells.Item(107, 14).Value = -4911

$ws.Cells.Item(134, 8).Value = 7805.4736
$ws.Cells.Item(134, 9).Value = 8356.5
$ws.Cells.Item(134, 11).Value = 25069.5
$ws.Cells.Item(134, 13).Value = -22534.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1881355.6
$ws.Cells.Item(31, 9).Value = 2858372.2
$ws.Cells.Item(31, 11).Value = 2858372.2
$ws.Cells.Item(31, 13).Value = -2858077.2

$ws.Cells.Item(34, 8).Value = 1881355.6
$ws.Cells.Item(34, 9).Value = 2858372.2
$ws.Cells.Item(34, 11).Value = 2858372.2
$ws.Cells.Item(34, 13).Value = -2858170.2

$ws.Cells.Item(107, 8).Value = 2855
$ws.Cells.Item(107, 9).Value = 711
$ws.Cells.Item(107, 11).Value = 711
$ws.Cells.Item(107, 13).Value = 1209

$ws.Cells.Item(122, 8).Value = 4235.5884
$ws.Cells.Item(122, 9).Value = 3396.75
$ws.Cells.Item(122, 11).Value = 10190.25
$ws.Cells.Item(122, 13).Value = -7740.25

$ws.Cells.Item(132, 8).Value = 1290.3334
$ws.Cells.Item(132, 9).Value = 1002.451
$ws.Cells.Item(132, 11).Value = 3007.353
$ws.Cells.Item(132, 13).Value = -477.3530000000001

$ws.Cells.Item(134, 8).Value = 1520.849
$ws.Cells.Item(134, 9).Value = 1346.5652
$ws.Cells.Item(134, 11).Value = 4039.6956
$ws.Cells.Item(134, 13).Value = -1504.6956

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 2523.25
$ws.Cells.Item(103, 9).Value = 2407.8333
$ws.Cells.Item(103, 10).Value = 2638.6667
$ws.Cells.Item(103, 11).Value = 7223.499899999999
$ws.Cells.Item(103, 12).Value = 7916.000100000001
$ws.Cells.Item(103, 13).Value = -6344.499899999999
$ws.Cells.Item(103, 14).Value = -9674.000100000001

$ws.Cells.Item(119, 8).Value = 50001104
$ws.Cells.Item(119, 9).Value = 62500884
$ws.Cells.Item(119, 11).Value = 187502652
$ws.Cells.Item(119, 13).Value = -187497814

$ws.Cells.Item(122, 8).Value = 912.8461
$ws.Cells.Item(122, 10).Value = 1129.52
$ws.Cells.Item(122, 12).Value = 10165.68
$ws.Cells.Item(122, 14).Value = -15065.68

$ws.Cells.Item(129, 8).Value = 52968.215
$ws.Cells.Item(129, 10).Value = 73785.8
$ws.Cells.Item(129, 12).Value = 221357.4
$ws.Cells.Item(129, 14).Value = -231357.4

$ws.Cells.Item(131, 8).Value = 5504311.5
$ws.Cells.Item(131, 10).Value = 10458.094
$ws.Cells.Item(131, 12).Value = 31374.282
$ws.Cells.Item(131, 14).Value = -41454.282

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2120.3704
$ws.Cells.Item(102, 9).Value = 2018.4
$ws.Cells.Item(102, 11).Value = 2018.4
$ws.Cells.Item(102, 13).Value = -396.4000000000001

$ws.Cells.Item(132, 8).Value = 1674860.9
$ws.Cells.Item(132, 9).Value = 2749357.8
$ws.Cells.Item(132, 11).Value = 8248073.399999999
$ws.Cells.Item(132, 13).Value = -8245543.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 10290.3
$ws.Cells.Item(16, 10).Value = 519.3333
$ws.Cells.Item(16, 12).Value = 519.3333
$ws.Cells.Item(16, 14).Value = -859.3333

$ws.Cells.Item(61, 8).Value = 2194.9583
$ws.Cells.Item(61, 9).Value = 1926.8334
$ws.Cells.Item(61, 10).Value = 2999.3333
$ws.Cells.Item(61, 11).Value = 1926.8334
$ws.Cells.Item(61, 12).Value = 2999.3333
$ws.Cells.Item(61, 13).Value = -1724.8334
$ws.Cells.Item(61, 14).Value = -3403.3333

$ws.Cells.Item(93, 8).Value = 1856.4
$ws.Cells.Item(93, 10).Value = 2021.75
$ws.Cells.Item(93, 12).Value = 2021.75
$ws.Cells.Item(93, 14).Value = -4517.75

$ws.Cells.Item(113, 8).Value = 2194.9583
$ws.Cells.Item(113, 9).Value = 1926.8334
$ws.Cells.Item(113, 10).Value = 2999.3333
$ws.Cells.Item(113, 11).Value = 1926.8334
$ws.Cells.Item(113, 12).Value = 2999.3333
$ws.Cells.Item(113, 13).Value = 243.1666
$ws.Cells.Item(113, 14).Value = -7339.3333

$ws.Cells.Item(128, 8).Value = 45000
$ws.Cells.Item(128, 10).Value = 45000
$ws.Cells.Item(128, 12).Value = 45000
$ws.Cells.Item(128, 14).Value = -54960

$ws.Cells.Item(132, 8).Value = 3362.7334
$ws.Cells.Item(132, 9).Value = 2650.5
$ws.Cells.Item(132, 10).Value = 4176.7144
$ws.Cells.Item(132, 11).Value = 7951.5
$ws.Cells.Item(132, 12).Value = 12530.1432
$ws.Cells.Item(132, 13).Value = -5421.5
$ws.Cells.Item(132, 14).Value = -17590.1432

$ws.Cells.Item(136, 8).Value = 2789.5557
$ws.Cells.Item(136, 9).Value = 1812.4286
$ws.Cells.Item(136, 11).Value = 5437.2858
$ws.Cells.Item(136, 13).Value = -2887.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 13984.333
$ws.Cells.Item(41, 10).Value = 13984.333
$ws.Cells.Item(41, 12).Value = 13984.333
$ws.Cells.Item(41, 14).Value = -14764.333

$ws.Cells.Item(45, 8).Value = 63712.332
$ws.Cells.Item(45, 9).Value = 20568.5
$ws.Cells.Item(45, 10).Value = 150000
$ws.Cells.Item(45, 11).Value = 20568.5
$ws.Cells.Item(45, 12).Value = 150000
$ws.Cells.Item(45, 13).Value = -20077.5
$ws.Cells.Item(45, 14).Value = -150982

$ws.Cells.Item(113, 8).Value = 400
$ws.Cells.Item(113, 10).Value = 600.75
$ws.Cells.Item(113, 12).Value = 1802.25
$ws.Cells.Item(113, 14).Value = -6142.25

$ws.Cells.Item(132, 8).Value = 1912.4642
$ws.Cells.Item(132, 9).Value = 1588.625
$ws.Cells.Item(132, 10).Value = 2344.25
$ws.Cells.Item(132, 11).Value = 4765.875
$ws.Cells.Item(132, 12).Value = 7032.75
$ws.Cells.Item(132, 13).Value = -2235.875
$ws.Cells.Item(132, 14).Value = -12092.75

$ws.Cells.Item(133, 8).Value = 79857.5
$ws.Cells.Item(133, 10).Value = 79857.5
$ws.Cells.Item(133, 12).Value = 79857.5
$ws.Cells.Item(133, 14).Value = -89977.5

$ws.Cells.Item(137, 8).Value = 38373.75
$ws.Cells.Item(137, 10).Value = 38373.75
$ws.Cells.Item(137, 12).Value = 38373.75
$ws.Cells.Item(137, 14).Value = -48573.75
